# TestRunnerBDD.xlsx - "commit after long time"
#  - FeatureSelection: best-fit width on column A
#  - Tag Selection: add a new "@firstflow" / "Yes" row, best-fit column A,
#    move the selection, set page orientation to portrait

$wb = $excel.ActiveWorkbook

$wsFeature = $wb.Worksheets.Item("FeatureSelection")
$wsTags    = $wb.Worksheets.Item("Tag Selection")

# --- FeatureSelection: column A best-fit (content already present) ---
$wsFeature.Columns.Item(1).ColumnWidth = 14.15

# --- Tag Selection: new data row ---
$wsTags.Range("A2").Value = "@firstflow"
$wsTags.Range("B2").Value = "Yes"

# Column A best-fit after the new, longer value was entered
$wsTags.Columns.Item(1).ColumnWidth = 11.3

# Page setup tweak recorded for this sheet
$wsTags.PageSetup.Orientation = 1

# Move the active selection (as last left by the author)
[void]$wsTags.Range("O16").Select()
